$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 53

# Copy the header style (bold, centered, bordered) from the existing
# last header cell (AC1) onto the three new header cells, then set
# their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row gets the same team record: 72 wins, 89 losses, 0 ties.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
